$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the style of the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
